# Data: - Formato matrix de datos - Matrix de datos del IDEAM por cuenca y año-mes
#
# The "detalle" sheet header row gets 3 new leading data columns
# (area, pend_md, zona_climatica) inserted right after "fecha" (col C),
# pushing the existing columns (elevacion..zona_climatica) to the right.
# The values for the three new headers are the ones that used to live at
# the tail of the table (old area/pend_md/zona_climatica columns), and the
# lone data cell in row 2 ("mediana") slides along with its column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("detalle")
$ws.Activate()

# Yellow fill matching the workbook's existing style for this class of
# header (same fill already used by def_acum/def_area/pend_md, etc.).
$amarillo = 65535   # FFFFFF00

# Insert 3 blank columns right before column D ("elevacion"), shifting
# elevacion..zona_climatica (D:R) to G:U.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header cells (D1:F1) - values move in from their old tail position
# (now R1, S1, U1 after the shift).
$ws.Range("D1").Value2 = "area"
$ws.Range("E1").Value2 = "pend_md"
$ws.Range("F1").Value2 = "zona_climatica"

# Remove the now-duplicated old header cells (their values now live in
# D1:F1). Delete right-to-left so the column letters stay valid while
# q_delta (T1) slides back in to become the new last column (R1).
$ws.Range("U1").EntireColumn.Delete()
$ws.Range("S1").EntireColumn.Delete()
$ws.Range("R1").EntireColumn.Delete()

# Give the new headers the same yellow fill used elsewhere for this kind
# of derived/basin metadata column.
$ws.Range("D1:F1").Interior.Color = $amarillo

# Update the saved cursor position to match.
$ws.Range("K9").Select()
